$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'72.741.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.12%  "
$ws.Range("D3").Value = "'3.977.10"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.48%  "
$ws.Range("D4").Value = "'1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'588.26"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +8.36%  "
$ws.Range("D6").Value = "'158.11"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.28%  "
$ws.Range("E7").Value = "  -3.38%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "'0.748"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("D10").Value = "'0.167"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.27%  "
$ws.Range("D11").Value = "'53.62"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.91%  "
$ws.Range("D12").Value = "'0.0000317"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.10%  "
$ws.Range("D13").Value = "'10.84"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.56%  "
$ws.Range("D14").Value = "'4.627.07"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.08%  "
$ws.Range("D15").Value = "'4.009.52"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.56%  "
$ws.Range("D16").Value = "'1.27"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +7.58%  "
$ws.Range("D17").Value = "'14.00"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.67%  "
$ws.Range("D18").Value = "'20.40"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.39%  "
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("D20").Value = "'72.604.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("D21").Value = "'430.72"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.13%  "
$ws.Range("E22").Value = "  +9.97%  "
$ws.Range("D23").Value = "'95.86"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.42%  "
$ws.Range("D24").Value = "'3.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -3.13%  "
$ws.Range("D25").Value = "'14.29"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.38%  "
$ws.Range("D26").Value = "'4.41"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +22.47%  "
$ws.Range("D27").Value = "'11.35"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.43%  "
$ws.Range("D28").Value = "'10.81"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.75%  "
$ws.Range("D29").Value = "'5.93"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +1.21%  "
$ws.Range("D30").Value = "'36.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.14%  "
$ws.Range("E31").Value = "  +7.48%  "
$ws.Range("D32").Value = "'50.88"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.61%  "
$ws.Range("E33").Value = "  +0.53%  "
$ws.Range("E34").Value = "  +0.03%  "
$ws.Range("D35").Value = "'679.10"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.38%  "
$ws.Range("D36").Value = "'68.36"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +3.62%  "
$ws.Range("D37").Value = "'0.438"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.91%  "
$ws.Range("E38").Value = "  +3.12%  "
$ws.Range("D39").Value = "'3.39"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.68%  "
$ws.Range("E40").Value = "  +0.12%  "
$ws.Range("E41").Value = "  -4.30%  "
$ws.Range("E42").Value = "  -3.00%  "
$ws.Range("B43").Value = "FirstDigitalUSD"
$ws.Range("C43").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D43").Value = "'1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.17%  "
$ws.Range("B44").Value = "THORChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D44").Value = "'10.95"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +11.68%  "
$ws.Range("D45").Value = "'0.0485"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.49%  "
$ws.Range("B46").Value = "Stellar"
$ws.Range("C46").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D46").Value = "'0.149"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.71%  "
$ws.Range("B47").Value = "Fetch.AI"
$ws.Range("C47").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D47").Value = "'2.69"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.67%  "
$ws.Range("D48").Value = "'3.40"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.88%  "
$ws.Range("D49").Value = "'3.42"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.02%  "
$ws.Range("D50").Value = "'2.99"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.26%  "
$ws.Range("E51").Value = "  +7.09%  "
